$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper used below to write plain text into a cell without Excel's
# automatic "this looks like a date" conversion (and without leaving a
# stray quote-prefix style behind): compute the text via a formula in a
# scratch cell, then paste-special just the value into the destination.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# Row 3 ("2019-05-15" / leche con avena / Ensalada de atún) closes out -
# only the leading id flag in column A stays.
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("F3").ClearContents()

# New row 4: next day's history entry.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 0

Set-TextValue $ws.Range("B4") "2019-05-16"
$ws.Range("C4").Value = 1
Set-TextValue $ws.Range("D4") "leche de soja con copos de maiz"
Set-TextValue $ws.Range("E4") "leche con avena"
